$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "57.005.43"
$ws.Range("E2").Value2 = "  -1.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.081.05"
$ws.Range("E3").Value2 = "  -0.54%  "
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "519.14"
$ws.Range("E5").Value2 = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "135.59"
$ws.Range("E6").Value2 = "  -3.85%  "
$ws.Range("E7").Value2 = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.082.17"
$ws.Range("E8").Value2 = "  -0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.456"
$ws.Range("E9").Value2 = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.34"
$ws.Range("E11").Value2 = "  -1.82%  "
$ws.Range("E12").Value2 = "  +1.82%  "
$ws.Range("E13").Value2 = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "3.612.05"
$ws.Range("E14").Value2 = "  -0.45%  "
$ws.Range("E15").Value2 = "  -1.64%  "
$ws.Range("E16").Value2 = "  -2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "57.125.95"
$ws.Range("E17").Value2 = "  -1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.084.02"
$ws.Range("E18").Value2 = "  -0.39%  "
$ws.Range("E19").Value2 = "  -3.72%  "
$ws.Range("E20").Value2 = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.83"
$ws.Range("E21").Value2 = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "347.19"
$ws.Range("E22").Value2 = "  +1.18%  "
$ws.Range("B23").Value2 = "Dai"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.00"
$ws.Range("E23").Value2 = "  +0.06%  "
$ws.Range("B24").Value2 = "LEO"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "5.77"
$ws.Range("E24").Value2 = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "68.26"
$ws.Range("E25").Value2 = "  +1.26%  "
$ws.Range("E26").Value2 = "  -2.90%  "
$ws.Range("E27").Value2 = "  -2.42%  "
$ws.Range("E28").Value2 = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.0₃0860"
$ws.Range("E29").Value2 = "  -6.46%  "
$ws.Range("E30").Value2 = "  -0.05%  "
$ws.Range("E31").Value2 = "  -0.32%  "
$ws.Range("E32").Value2 = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "5.80"
$ws.Range("E33").Value2 = "  -9.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "20.80"
$ws.Range("E34").Value2 = "  -0.88%  "
$ws.Range("E35").Value2 = "  +5.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "159.42"
$ws.Range("E36").Value2 = "  +0.60%  "
$ws.Range("E37").Value2 = "  -4.34%  "
$ws.Range("E38").Value2 = "  -2.75%  "
$ws.Range("E39").Value2 = "  -2.36%  "
$ws.Range("E40").Value2 = "  -1.09%  "
$ws.Range("E41").Value2 = "  -2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "4.05"
$ws.Range("E42").Value2 = "  +0.46%  "
$ws.Range("E43").Value2 = "  +1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.691"
$ws.Range("E44").Value2 = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "2.379.95"
$ws.Range("E45").Value2 = "  +4.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "36.59"
$ws.Range("E46").Value2 = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "3.120.61"
$ws.Range("E48").Value2 = "  -0.50%  "
$ws.Range("E49").Value2 = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.951"
$ws.Range("E50").Value2 = "  -4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "5.94"
$ws.Range("E51").Value2 = "  -2.73%  "
